$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report Rcp")
$ws2 = $wb.Worksheets.Item("Report Sales")

# Remember the existing comments (cell -> text) before we shift rows, so we
# can re-create them one row lower afterwards (the engine's row-Insert does
# not itself relocate cell comments).
$commentMap = @{}
foreach ($addr in @("C3", "C5", "C6", "F3", "F4", "F5")) {
    $cmt = $ws.Range($addr).Comment
    if ($cmt -ne $null) {
        $commentMap[$addr] = $cmt.Text()
    }
}
foreach ($addr in $commentMap.Keys) {
    $ws.Range($addr).Comment.Delete()
}

# Insert a new row above row 2: row 1 (the report title) keeps its
# formatting but row 2 becomes a duplicate of it, and every row from the
# old row 3 onward shifts down by one.
$ws.Rows("2:2").Insert()

# Move the title text out of row 1 and into the newly inserted row 2,
# restoring the merged header cell on its new row.
$ws.Range("B1").Cut($ws.Range("B2"))
$ws.Range("B1:F1").UnMerge()
$ws.Range("B2:F2").Merge()

# Re-create the comments one row below their original position.
$ws.Range("C4").AddComment($commentMap["C3"]) | Out-Null
$ws.Range("C6").AddComment($commentMap["C5"]) | Out-Null
$ws.Range("C7").AddComment($commentMap["C6"]) | Out-Null
$ws.Range("F4").AddComment($commentMap["F3"]) | Out-Null
$ws.Range("F5").AddComment($commentMap["F4"]) | Out-Null
$ws.Range("F6").AddComment($commentMap["F5"]) | Out-Null

# Leave the whole first row selected on both sheets, matching the
# post-edit selection state.
$ws2.Rows("1:1").Select()
$ws.Rows("1:1").Select()
